# Scheduled-runner style update of the per-sheet price/profit columns
# (H..N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, reflecting
# refreshed market-board pricing data. Cells whose computed profit is no
# longer meaningful for a row are cleared to match the source data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 346.83334
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -328

$ws.Range("H5").Value = 87.25
$ws.Range("I5").Value = 87.25
$ws.Range("K5").Value = 87.25
$ws.Range("M5").Value = 27.75

$ws.Range("H6").Value = 519.7143
$ws.Range("I6").Value = 435
$ws.Range("J6").Value = 583.25
$ws.Range("K6").Value = 1305
$ws.Range("L6").Value = 1749.75
$ws.Range("M6").Value = -1193
$ws.Range("N6").Value = -1973.75

$ws.Range("H40").Value = 2266.6667
$ws.Range("J40").Value = 1800
$ws.Range("L40").Value = 1800
$ws.Range("N40").Value = -2150

$ws.Range("H58").Value = 1907.5
$ws.Range("J58").Value = 1951.4286
$ws.Range("L58").Value = 5854.2858
$ws.Range("N58").Value = -6154.2858

$ws.Range("H92").Value = 325
$ws.Range("J92").Value = 325
$ws.Range("L92").Value = 325
$ws.Range("N92").Value = -2821

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws.Range("H138").Value = 1643.8
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11064.333
$ws.Range("I32").Value = 11064.333
$ws.Range("K32").Value = 11064.333
$ws.Range("M32").Value = -10777.333

$ws.Range("H61").Value = 3598.1428
$ws.Range("I61").Value = 3598.1428
$ws.Range("K61").Value = 3598.1428
$ws.Range("M61").Value = -3386.1428

$ws.Range("H74").Value = 3599.6667
$ws.Range("I74").Value = 2899.5
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2899.5
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -2025.5
$ws.Range("N74").Value = -6748

$ws.Range("H76").Value = 55833.168
$ws.Range("J76").Value = 55833.168
$ws.Range("L76").Value = 55833.168
$ws.Range("N76").Value = -56509.168

$ws.Range("H77").Value = 3599.6667
$ws.Range("I77").Value = 2899.5
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 14497.5
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -10129.5
$ws.Range("N77").Value = -33736

$ws.Range("H79").Value = 55833.168
$ws.Range("J79").Value = 55833.168
$ws.Range("L79").Value = 55833.168
$ws.Range("N79").Value = -58173.168

$ws.Range("H132").Value = 2228.125
$ws.Range("I132").Value = 1710
$ws.Range("K132").Value = 5130
$ws.Range("M132").Value = -2600

$ws.Range("H136").Value = 3598.1428
$ws.Range("I136").Value = 3598.1428
$ws.Range("K136").Value = 10794.4284
$ws.Range("M136").Value = -8244.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1635.3846
$ws.Range("I134").Value = 1637.1666
$ws.Range("K134").Value = 4911.4998
$ws.Range("M134").Value = -2376.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -5278

$ws.Range("H22").Value = 401.16666
$ws.Range("I22").Value = 401.16666
$ws.Range("K22").Value = 401.16666
$ws.Range("M22").Value = -51.16665999999998

$ws.Range("H31").Value = 2964.1
$ws.Range("I31").Value = 1761.2
$ws.Range("K31").Value = 1761.2
$ws.Range("M31").Value = -1466.2

$ws.Range("H34").Value = 2964.1
$ws.Range("I34").Value = 1761.2
$ws.Range("K34").Value = 1761.2
$ws.Range("M34").Value = -1559.2

$ws.Range("H35").Value = 12055
$ws.Range("I35").Value = 8749
$ws.Range("J35").Value = 17014
$ws.Range("K35").Value = 8749
$ws.Range("L35").Value = 17014
$ws.Range("M35").Value = -8455
$ws.Range("N35").Value = -17602

$ws.Range("H58").Value = 2318.3333
$ws.Range("I58").Value = 1004
$ws.Range("K58").Value = 1004
$ws.Range("M58").Value = -801

$ws.Range("H96").Value = 24204.4
$ws.Range("J96").Value = 24204.4
$ws.Range("L96").Value = 24204.4
$ws.Range("N96").Value = -29696.4

$ws.Range("H107").Value = 810.3333
$ws.Range("I107").Value = 835
$ws.Range("J107").Value = 761
$ws.Range("K107").Value = 835
$ws.Range("L107").Value = 761
$ws.Range("M107").Value = 1085
$ws.Range("N107").Value = -4601

$ws.Range("H123").Value = 150999
$ws.Range("I123").Value = 150999
$ws.Range("K123").Value = 150999
$ws.Range("M123").Value = -146099

$ws.Range("H136").Value = 2318.3333
$ws.Range("I136").Value = 1004
$ws.Range("K136").Value = 3012
$ws.Range("M136").Value = -462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H38").Value = 1008.5
$ws.Range("J38").Value = 91.40000000000001
$ws.Range("L38").Value = 274.2
$ws.Range("N38").Value = -968.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3669483.2
$ws.Range("I11").Value = 3336133.8
$ws.Range("K11").Value = 3336133.8
$ws.Range("M11").Value = -3335994.8

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5586.7334
$ws.Range("I16").Value = 4221.5713
$ws.Range("K16").Value = 4221.5713
$ws.Range("M16").Value = -4051.5713

$ws.Range("H46").Value = 1137.5
$ws.Range("I46").Value = 1183.3334
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1183.3334
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -995.3334
$ws.Range("N46").Value = -1376

$ws.Range("H136").Value = 100497.2
$ws.Range("I136").Value = 5634.25
$ws.Range("K136").Value = 16902.75
$ws.Range("M136").Value = -14352.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 36790
$ws.Range("J69").Value = 36790
$ws.Range("L69").Value = 36790
$ws.Range("N69").Value = -38288

$ws.Range("H72").Value = 36790
$ws.Range("J72").Value = 36790
$ws.Range("L72").Value = 110370
$ws.Range("N72").Value = -117858

$ws.Range("H96").Value = 1635.8
$ws.Range("I96").Value = 1570
$ws.Range("J96").Value = 1734.5
$ws.Range("K96").Value = 1570
$ws.Range("L96").Value = 1734.5
$ws.Range("M96").Value = -197
$ws.Range("N96").Value = -4480.5

$ws.Range("H136").Value = 2784.6
$ws.Range("I136").Value = 2784.6
$ws.Range("K136").Value = 8353.799999999999
$ws.Range("M136").Value = -5803.799999999999
